# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
# Swap the data rows for the two workers currently on rows 16 and 17
# (C:G) so that row 16 now holds EINER ALIPIO MORALES MARTINEZ's record
# and row 17 holds MANUEL IVAN MORALES ARNEDO's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$cols = @("C", "D", "E", "F", "G")

$row16vals = @{}
$row17vals = @{}

foreach ($col in $cols) {
    $row16vals[$col] = $ws.Range("$col" + "16").Value()
    $row17vals[$col] = $ws.Range("$col" + "17").Value()
}

foreach ($col in $cols) {
    $ws.Range("$col" + "16").Value = $row17vals[$col]
    $ws.Range("$col" + "17").Value = $row16vals[$col]
}
